# Fix typo in control table: the DMX Rose control-code block for rows
# 25-30 and 31-36 was a duplicate of the codes used in rows 13-18 / 19-24.
# They should instead continue the numbering sequence (70-75 / 80-85).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$fixes = @{
    "O25" = 70
    "O26" = 71
    "O27" = 72
    "O28" = 73
    "O29" = 74
    "O30" = 75
    "O31" = 80
    "O32" = 81
    "O33" = 82
    "O34" = 83
    "O35" = 84
    "O36" = 85
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# Match the author's final selection (the corrected table range).
$excel.Goto($ws.Range("O2:S36"))
